# Apply the "Model Stub" -> "SAPL.io" tags.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# 1. Rename the worksheet (this also drives the workbook.xml <sheet name="..."> entry)
$ws.Name = "SAPL.io"

# 2. Extend formatting ranges to cover the new, bigger matrix (7 tag columns x 14 data rows)
#    before we touch any values, so every cell keeps the right style.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B1:G1").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2:A15").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B2:G15").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# 3. Header row - tag names across columns B:G
$ws.Range("A1").Value = "Element"
$ws.Range("B1").Value = "sapl-jwt"
$ws.Range("C1").Value = "sapl-spring-pdp-embedded"
$ws.Range("D1").Value = "sapl-webflux-endpoint"
$ws.Range("E1").Value = "spring-boot"
$ws.Range("F1").Value = "spring-security"
$ws.Range("G1").Value = "windows-server"

# 4. Element (row label) column, rows 2-15
$ws.Range("A2").Value  = "Customer Web Client"
$ws.Range("A3").Value  = "Customer Traffic"
$ws.Range("A4").Value  = "Policy Fileserver"
$ws.Range("A5").Value  = "SAPL MVC"
$ws.Range("A6").Value  = "SAPL Server LT Traffic"
$ws.Range("A7").Value  = "SAPL Server LT"
$ws.Range("A8").Value  = "Policy Filesystem Access"
$ws.Range("A9").Value  = "Access-Decision"
$ws.Range("A10").Value = "Access-Request"
$ws.Range("A11").Value = "Patient Data"
$ws.Range("A12").Value = "Policies"
$ws.Range("A13").Value = "File-Server Network"
$ws.Range("A14").Value = "SAPL-Server DMZ"
$ws.Range("A15").Value = "Web DMZ"

# 5. Clear out the whole matrix body first (keeps per-cell style, removes stale values/marks)
$ws.Range("B2:G15").ClearContents()

# 6. Re-apply the "X" marks at their new matrix positions
$ws.Range("G4").Value = "X"
$ws.Range("E5").Value = "X"
$ws.Range("B7").Value = "X"
$ws.Range("C7").Value = "X"
$ws.Range("D7").Value = "X"
$ws.Range("E7").Value = "X"
$ws.Range("F7").Value = "X"

# 7. Widen the tag columns (B:G) to match column A's "customWidth" 35 formatting
$ws.Range("B1:G1").EntireColumn.ColumnWidth = 34.16666666666666
